$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Updated simulation results for rows 2-8 (d=1..d=10), columns B-E
$values = @{
    2 = @(39.42036321020625, 39.2783745506574,  41.62977214893819, 42.55154529513415)
    3 = @(50.38248768211463, 50.31006197060847, 50.2463590208683,  50.28366705008921)
    4 = @(98.78855299310007, 98.72121748546154, 98.80617458215927, 98.89710136385391)
    5 = @(98.96221427295416, 98.96041453209654, 98.892789947442,   98.90607321957347)
    6 = @(98.55406822997308, 98.48081443726426, 98.49698157629641, 98.47964947405534)
    7 = @(97.59537080693012, 97.54753994293857, 97.64033857574367, 97.57654125726216)
    8 = @(96.334275508866,   96.08004943890383, 96.27314898098837, 96.15330927021672)
}

foreach ($row in $values.Keys) {
    $rowVals = $values[$row]
    $ws.Cells.Item($row, 2).Value = $rowVals[0]
    $ws.Cells.Item($row, 3).Value = $rowVals[1]
    $ws.Cells.Item($row, 4).Value = $rowVals[2]
    $ws.Cells.Item($row, 5).Value = $rowVals[3]
}
